$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (AddCustomerTest): drop the alerttext column (D), replace the
#     sample customer row with a new data-provider row ---
$ws1.Range("D1:D2").ClearContents()

$ws1.Range("A2").Value = "dang"
$ws1.Range("B2").Value = "chau"
$ws1.Range("C2").Value = 2

$ws1.Columns.Item(4).ColumnWidth = 27.140625

# --- Sheet2 (OpenAccountTest): new customer/currency data row ---
$ws2.Range("A2").Value = "dang chau"
$ws2.Range("B2").Value = "Dollar"

$ws2.Columns.Item(1).ColumnWidth = 9.85546875
$ws2.Columns.Item(2).ColumnWidth = 8.5703125

# Move the selection on sheet2 (this briefly activates it), then restore
# sheet1 as the active/selected tab to match the saved view state.
[void]$ws2.Range("H14").Select()
[void]$ws1.Select()
[void]$ws1.Range("D7").Select()
